# Update column F (dSF) values for specific rows to match the repulled data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = -2
$ws.Range("F5").Value = 11
$ws.Range("F8").Value = -2
$ws.Range("F9").Value = -4
$ws.Range("F16").Value = 0
$ws.Range("F18").Value = -1
$ws.Range("F22").Value = -5
$ws.Range("F23").Value = 9
$ws.Range("F24").Value = 2
